$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows 344-357 (columns A: date serial, B: nuovi pos., C: somma mobile 7gg., D: somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @(44418, 1, 6, 68.99724011039559),
    @(44419, 0, 6, 68.99724011039559),
    @(44420, 0, 5, 57.49770009199631),
    @(44421, 0, 5, 57.49770009199631),
    @(44422, 1, 5, 57.49770009199631),
    @(44423, 0, 3, 34.4986200551978),
    @(44424, 3, 5, 57.49770009199631),
    @(44425, 1, 5, 57.49770009199631),
    @(44426, 0, 5, 57.49770009199631),
    @(44427, 1, 6, 68.99724011039559),
    @(44428, 0, 6, 68.99724011039559),
    @(44429, 2, 7, 80.49678012879485),
    @(44430, 2, 9, 103.4958601655934),
    @(44431, 0, 6, 68.99724011039559)
)

$startRow = 344
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $prevRow = $r - 1

    # Copy the formatting (style) of the cell directly above, so new cells
    # reuse the existing "date" cell style instead of creating a new one.
    $ws.Cells.Item($prevRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$excel.CutCopyMode = 0
